# Updates cryptos list values (price/volume) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.081.67"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.644.56"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'217.51"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'0.5201"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.2617"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "'0.06281"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'20.41"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("D11").Value = "'0.07757"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "'4.477"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "1.634.67"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").Value = "1.870.92"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "'0.5579"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "0.0₅7991"
$ws.Range("D17").Value = "'64.78"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "26.077.16"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'4.642"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'192.44"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'10.10"
$ws.Range("D23").Value = "'5.950"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("D25").Value = "'146.29"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'0.1202"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "'7.165"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "'15.94"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "'1.485"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'0.05619"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("E32").Value = "  -4.25%  "
$ws.Range("D33").Value = "'3.368"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'2.788"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D37").Value = "'0.9366"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("D38").Value = "'0.5661"
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("D39").Value = "'5.955"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "'0.01581"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "1.052.15"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").Value = "'2.568"
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'0.8417"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'102.18"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("D46").Value = "1.782.13"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "'57.12"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "'1.009"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'0.05320"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4327"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.911"
$ws.Range("E51").Value = "  -1.23%  "
